# Insert one new weekly price record for "Apio" (Femacal de La Calera, Coquimbo)
# immediately before the existing row 720, shifting all subsequent rows down by
# one (old row 720 becomes 721, ..., old row 759 becomes row 760).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 720:759 down to 721:760, leaving a blank row 720 (inherits the
# formatting of the row being pushed down, same as a manual Excel row insert).
$ws.Rows("720:720").Insert()

# Populate the newly inserted row 720 with the new record's data.
$ws.Cells.Item(720, 1).Value  = 3
$ws.Cells.Item(720, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(720, 3).Value  = "Coquimbo"
$ws.Cells.Item(720, 4).Value  = 45267
$ws.Cells.Item(720, 5).Value  = 5
$ws.Cells.Item(720, 6).Value  = 100112017
$ws.Cells.Item(720, 7).Value  = "Apio"
$ws.Cells.Item(720, 8).Value  = "Americana (o)"
$ws.Cells.Item(720, 9).Value  = "Primera"
$ws.Cells.Item(720, 10).Value = 200
$ws.Cells.Item(720, 11).Value = 9000
$ws.Cells.Item(720, 12).Value = 9500
$ws.Cells.Item(720, 13).Value = 9250
$ws.Cells.Item(720, 14).Value = "`$/docena de matas"
$ws.Cells.Item(720, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(720, 16).Value = 1542
$ws.Cells.Item(720, 17).Value = 6
$ws.Cells.Item(720, 18).Value = "Hortaliza"
